# IR remote codes.xlsx -- "basically working. fade control still need work."
#
# 1. Rename Sheet1 -> "black elgo"
# 2. Add a new sheet "Snapper" after it
# 3. Rework "black elgo" column C/D: C keeps the MAP labels (reindexed after a
#    removed label), D switches from a duplicate-hex-string column to the
#    actual decimal IR int values
# 4. Populate "Snapper" with the new LED-strip remote mapping
# 5. Restore the selections recorded in each sheet's view

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "black elgo"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Snapper"

# ---------------------------------------------------------------------
# "black elgo" (sheet1) -- string cells, written in original edit order
# ---------------------------------------------------------------------
$ws1.Range("C1").Value = "MAP"
$ws1.Range("C14").Value = "red up"
$ws1.Range("C17").Value = "red down"
$ws1.Range("C15").Value = "green up"
$ws1.Range("C16").Value = "blue up"
$ws1.Range("C18").Value = "green down"
$ws1.Range("C19").Value = "blue down"
$ws1.Range("C4").Value = "all off"
$ws1.Range("C2").Value = "all on"

# ---------------------------------------------------------------------
# "Snapper" (sheet2) -- string cells, written in original edit order
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = "bright up"
$ws2.Range("C3").Value = "bright down"
$ws2.Range("C4").Value = "off"
$ws2.Range("C5").Value = "on"
$ws2.Range("C11").Value = "red"
$ws2.Range("C9").Value = "blue"
$ws2.Range("C10").Value = "green"
$ws2.Range("C8").Value = "white"
$ws2.Range("C12").Value = "flash"
$ws2.Range("D1").Value = "row "
$ws2.Range("C15").Value = "strobe"
$ws2.Range("C19").Value = "fade"
$ws2.Range("C22").Value = "smooth"
$ws2.Range("F7").Value = "Colors designate columns via the first row buttons"
$ws2.Range("F9").Value = "NB, several red column buttons do not work"

# Last of the new labels -- lands at the end of the shared-string table
$ws1.Range("D1").Value = "int values"

# ---------------------------------------------------------------------
# "black elgo" -- column D: decimal IR int values (replaces old hex dupes)
# ---------------------------------------------------------------------
$ws1.Range("D2").Value = 16753245
$ws1.Range("D3").Value = 16736925
$ws1.Range("D4").Value = 16769565
$ws1.Range("D5").Value = 16720605
$ws1.Range("D6").Value = 16712445
$ws1.Range("D7").Value = 16761405
$ws1.Range("D8").Value = 16769055
$ws1.Range("D9").Value = 16754775
$ws1.Range("D10").Value = 16748655
$ws1.Range("D11").Value = 16738455
$ws1.Range("D12").Value = 16750695
$ws1.Range("D13").Value = 16756815
$ws1.Range("D14").Value = 16724175
$ws1.Range("D15").Value = 16718055
$ws1.Range("D16").Value = 16743045
$ws1.Range("D17").Value = 16716015
$ws1.Range("D18").Value = 16726215
$ws1.Range("D19").Value = 16734885
$ws1.Range("D20").Value = 16728765
$ws1.Range("D21").Value = 16730805
$ws1.Range("D22").Value = 16732845

# ---------------------------------------------------------------------
# "Snapper" -- column B (raw IR codes) and D (row-group numbers)
# ---------------------------------------------------------------------
$ws2.Range("B2").Value = 16187647
$ws2.Range("D2").Value = 1
$ws2.Range("B3").Value = 16220287
$ws2.Range("D3").Value = 1
$ws2.Range("B4").Value = 16203967
$ws2.Range("D4").Value = 1
$ws2.Range("B5").Value = 16236599
$ws2.Range("D5").Value = 1
$ws2.Range("B6").Value = 16220287
$ws2.Range("B7").Value = 16187647
$ws2.Range("B8").Value = 16244759
$ws2.Range("D8").Value = 2
$ws2.Range("B9").Value = 16212127
$ws2.Range("D9").Value = 2
$ws2.Range("B10").Value = 16228447
$ws2.Range("D10").Value = 2
$ws2.Range("B11").Value = 16195807
$ws2.Range("D11").Value = 2
$ws2.Range("B12").Value = 16240675
$ws2.Range("D12").Value = 3
$ws2.Range("B13").Value = 16208047
$ws2.Range("C13").Value = "blue"
$ws2.Range("D13").Value = 3
$ws2.Range("B14").Value = 16224359
$ws2.Range("C14").Value = "green"
$ws2.Range("D14").Value = 3
$ws2.Range("B15").Value = 16248839
$ws2.Range("D15").Value = 4
$ws2.Range("B16").Value = 16248972
$ws2.Range("C16").Value = "blue"
$ws2.Range("D16").Value = 4
$ws2.Range("B17").Value = 16232527
$ws2.Range("C17").Value = "green"
$ws2.Range("D17").Value = 4
$ws2.Range("B18").Value = 16199887
$ws2.Range("C18").Value = "red"
$ws2.Range("D18").Value = 4
$ws2.Range("B19").Value = 16238647
$ws2.Range("D19").Value = 5
$ws2.Range("B20").Value = 16206007
$ws2.Range("C20").Value = "blue"
$ws2.Range("D20").Value = 5
$ws2.Range("B21").Value = 16222327
$ws2.Range("C21").Value = "green"
$ws2.Range("D21").Value = 5
$ws2.Range("B22").Value = 16246790
$ws2.Range("D22").Value = 6
$ws2.Range("B23").Value = 16246933
$ws2.Range("C23").Value = "blue"
$ws2.Range("D23").Value = 6
$ws2.Range("B24").Value = 16230470
$ws2.Range("C24").Value = "green"
$ws2.Range("D24").Value = 6
$ws2.Range("B25").Value = 16246918
$ws2.Range("C25").Value = "red"
$ws2.Range("D25").Value = 6

# ---------------------------------------------------------------------
# Selections recorded in the saved workbook
# ---------------------------------------------------------------------
$ws2.Range("F9").Select() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("C5").Select() | Out-Null
